$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The authoritative list of scheme names for column B, rows 2-19, in the
# order they must be (re-)interned into the shared-string table so that
# the table comes out in the same order as the target workbook:
#   HKL, ND Single, RD Single, TD Single, Morris,
#   Ring Perpendicular to ND/RD/TD,
#   Gaussian-Quadrature,
#   Spiral-90deg-10rot-5space, Spiral-90deg-15rot-5space, Spiral-90deg-10rot-3space,
#   NoRotation-tilt60deg, Rotation-NoTilt, Rotation-60detTilt,
#   HexGrid-90degTilt5degRes, HexGrid-90degTilt22p5degRes, HexGrid-60degTilt5degRes
$schemeNames = @(
    "HKL",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "Gaussian-Quadrature",
    "Spiral-90deg-10rot-5space",
    "Spiral-90deg-15rot-5space",
    "Spiral-90deg-10rot-3space",
    "NoRotation-tilt60deg",
    "Rotation-NoTilt",
    "Rotation-60detTilt",
    "HexGrid-90degTilt5degRes",
    "HexGrid-90degTilt22p5degRes",
    "HexGrid-60degTilt5degRes"
)

# Clear out column B's existing string content so the shared-string table
# gets rebuilt from scratch, in the write order below (this is what lets
# us re-order the shared strings to match the reran-script output).
$ws.Range("B1:B19").ClearContents()

for ($i = 0; $i -lt $schemeNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $schemeNames[$i]
}

# Extend the table down to the three new rows (17-19), duplicating the
# HexGrid rows the same way the rerun notebook did, with column A as the
# running 0-based index and C:M averaged-intensity placeholders of 1.
for ($row = 17; $row -le 19; $row++) {
    $idx = $row - 2

    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $idx
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    for ($col = 3; $col -le 13; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}
